{"js": "// Removed eff date from Fowler scheduling template.\n//\n// 1) \"...Judge Mark W. Fowler \u2013 effective March 25, 2025\" ->\n//    \"...Judge Mark W. Fowler \" (extra leading tab before \"Judge\", and the\n//    \" \u2013 effective March 25, 2025\" suffix dropped).\n// 2) Footer \"Printed\" date field result text updated to match the new\n//    effective date/time.\n\nconst body = context.document.body;\nconst target = body.search(\"Judge Mark W. Fowler \\u2013 effective March 25, 2025\", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\"\\tJudge Mark W. Fowler \", \"Replace\");\n  await context.sync();\n}\n\n// The footer's PRINTDATE field result lives in the first-page footer story.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const section of sections.items) {\n  for (const type of [\"FirstPage\", \"Primary\", \"EvenPages\"]) {\n    let footerBody;\n    try {\n      footerBody = section.getFooter(type);\n    } catch (e) {\n      continue;\n    }\n    const hits = footerBody.search(\"3/1/2024 8:33 AM\", { matchCase: true });\n    hits.load(\"text\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[0].insertText(\"3/25/2025 11:29 AM\", \"Replace\");\n      await context.sync();\n    }\n  }\n}\n", "ps1": "# Removed eff date from Fowler scheduling template.\n#\n# 1) \"...Judge Mark W. Fowler - effective March 25, 2025\" ->\n#    \"...Judge Mark W. Fowler \" (extra leading tab before \"Judge\", and the\n#    \" - effective March 25, 2025\" suffix dropped).\n# 2) Footer \"Printed\" date field result text updated to match the new\n#    effective date/time.\n\n$d = $word.ActiveDocument\n\n# --- 1) Main body: drop the effective-date suffix and add a tab. ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Judge Mark W. Fowler \" + [char]0x2013 + \" effective March 25, 2025\"\n$find.Replacement.Text = [char]9 + \"Judge Mark W. Fowler \"\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n\n# --- 2) Footer: refresh the cached PRINTDATE field-result text. ---\nforeach ($sec in $d.Sections) {\n    foreach ($idx in 1, 2, 3) {\n        $footer = $sec.Footers.Item($idx)\n        if ($footer.Exists -and $footer.Range.Text.Contains(\"3/1/2024 8:33 AM\")) {\n            $ffind = $footer.Range.Find\n            $ffind.ClearFormatting()\n            $ffind.Replacement.ClearFormatting()\n            $ffind.Text = \"3/1/2024 8:33 AM\"\n            $ffind.Replacement.Text = \"3/25/2025 11:29 AM\"\n            $ffind.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n        }\n    }\n}\n"}
